$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# 14. Nº Horas: 40 -> 42
# (done first: the replacement text for the R.U.T. field below ends up
#  containing the substring "40", which would otherwise get caught by
#  this Find if run afterwards)
Replace-Text "40" "42"

# 1. CAI N°5 -> CAI N°1
Replace-Text "CAI N°5" "CAI N°1"

# 2. C.abierto -> C.cerrado
Replace-Text "C.abierto" "C.cerrado"

# 3. FECHA 2022-01-07 -> 2022-01-09
Replace-Text "2022-01-07" "2022-01-09"

# 4. RAZÓN SOCIAL value: GERARD CANO GÓMEZ -> CENCOSUD S.A.
Replace-Text "GERARD CANO GÓMEZ" "CENCOSUD S.A."

# 5-7. Empty GIRO / ATENCIÓN / DEPARTAMENTO values -> "asd"
# (Find cannot match empty text, so address the table cells directly.
#  Table.Cell(r,c) mis-handles gridSpan in this table, so use Row.Cells.)
$tCliente = $d.Tables.Item(3)
$tCliente.Rows.Item(3).Cells.Item(2).Range.Text = "asd"   # GIRO:
$tCliente.Rows.Item(4).Cells.Item(2).Range.Text = "asd"   # ATENCIÓN
$tCliente.Rows.Item(5).Cells.Item(2).Range.Text = "asd"   # DEPARTAMENTO

# 8. R.U.T. value: 12076136-6 -> 12621140-6
Replace-Text "12076136-6" "12621140-6"

# 9. DIRECCION / COMUNA value: "/" -> "/maipu"
# ("/" is not unique as a Find target because it is also a substring of the
#  "DIRECCION / COMUNA" label, so set the cell directly.)
$tCliente.Rows.Item(7).Cells.Item(2).Range.Text = "/maipu"

# 10. FONO value: 936577225 -> 123123123
Replace-Text "936577225" "123123123"

# 11. VENCIMIENTO value: 02/10/2022 -> 21/12/2021
Replace-Text "02/10/2022" "21/12/2021"

# 12. Curso name
Replace-Text "Aplicación De Técnicas De Control De Inventarios" "Aplicación De Herramientas De Redacción Y Ortografía Para Profesionales"

# 13. Código Sence
Replace-Text "1238020234" "1238020242"

# 15. Valor Unitario / Valor Total / TOTAL: 160000 -> 168000 (x3, all identical)
Replace-Text "160000" "168000"

# 16-17. Checkbox "X" moves from USACH to EMPRESA
$tCheck = $d.Tables.Item(5)
$tCheck.Cell(1,1).Range.Text = ""
$tCheck.Cell(1,4).Range.Text = "X"

# 18. OTRO: None -> OTRO: (empty)
Replace-Text "None" ""

# 19. Orden de Compra N°12345 -> N°123123
Replace-Text "12345" "123123"

# 20. Obs: Prueba factura null -> Obs: ivan weco 2
Replace-Text "Prueba factura null" "ivan weco 2"
